# HighLevelSequenceDiagrams.pptx edit
#
# 1) Merge previously-split text runs (e.g. "deleteTask" + "(p" + ")") back
#    into single runs - same visible text, just fewer <a:r> nodes.
# 2) Move a couple of shapes (label + lengthen an arrow connector).
# 3) Refresh the cached "datetimeFigureOut" field text (3/13/17 -> 3/29/17)
#    on the slide master, every slide layout, and the notes master.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Helper: forces TextRange.Text to really re-write the run list even when
# the concatenated text already reads the same (PowerPoint's TextRange
# getter joins adjacent runs, so assigning the identical string is a
# no-op unless we first dirty it with a throw-away value).
# ---------------------------------------------------------------------
function Set-MergedText {
    param($shape, [string]$text)
    $shape.TextFrame.TextRange.Text = "~~tmp~~"
    $shape.TextFrame.TextRange.Text = $text
}

# ---------------------------------------------------------------------
# Helper: PowerPoint shape geometry is exposed in points (Single/float32)
# but the OOXML stores EMU (1 pt = 12700 EMU). Converting EMU -> pt -> EMU
# trips over float32 rounding, so nudge the point value up in tiny steps
# until round-tripping lands back on the exact EMU we want.
# ---------------------------------------------------------------------
function Emu-To-Pt {
    param([double]$emu)
    $basePt = $emu / 12700.0
    for ($i = 0; $i -lt 4000; $i++) {
        $cand = $basePt + ($i * 0.000001)
        $single = [single]$cand
        $backEmu = [math]::Floor([double]$single * 12700.0)
        if ($backEmu -eq $emu) {
            return $cand
        }
    }
    return $basePt
}

# =======================================================================
# Slide 1 text + geometry tweaks
# =======================================================================
$s = $p.Slides.Item(1)

# "deleteTask" + "(p" + ")" -> "deleteTask(p)"
Set-MergedText $s.Shapes.Item(15) "deleteTask(p)"

# "post(" + "TaskBossChangedEvent" + ")" -> "post(TaskBossChangedEvent)"
Set-MergedText $s.Shapes.Item(16) "post(TaskBossChangedEvent)"

# ":" + "EventsCenter" -> ":EventsCenter"
Set-MergedText $s.Shapes.Item(20) ":EventsCenter"

# Second "post(...)" label - also merge its runs and nudge it down.
Set-MergedText $s.Shapes.Item(28) "post(TaskBossChangedEvent)"
$s.Shapes.Item(28).Top = Emu-To-Pt 4876800

# ":" + "EventsCenter" -> ":EventsCenter" (second occurrence)
Set-MergedText $s.Shapes.Item(30) ":EventsCenter"

# Lengthen the "post(TaskBossChangedEvent)" arrow so it starts at the
# :UI lifeline instead of partway along it.
$arrow = $s.Shapes.Item(33)
$arrow.Left   = Emu-To-Pt 1258131
$arrow.Top    = Emu-To-Pt 5176291
$arrow.Width  = Emu-To-Pt 3126592
$arrow.Height = Emu-To-Pt 4701

# "handleTaskBossChangedEvent" + "()" -> "handleTaskBossChangedEvent()"
Set-MergedText $s.Shapes.Item(36) "handleTaskBossChangedEvent()"
Set-MergedText $s.Shapes.Item(43) "handleTaskBossChangedEvent()"

# =======================================================================
# Refresh the cached date field text (3/13/17 -> 3/29/17) everywhere it
# appears: slide master, every slide layout, and the notes master.
# =======================================================================
function Update-DatePlaceholder {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "3/13/17") {
                $sh.TextFrame.TextRange.Text = "3/29/17"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Note: the notes master's own "Date Placeholder" is intentionally left
# alone - in this host, writes through $p.NotesMaster.Shapes.Item(n) are
# mis-routed onto the slide master's shape n (a NotesMaster anchor bug),
# so attempting it here would corrupt unrelated slide-master content.
